$d = $word.ActiveDocument
$w = "http://schemas.openxmlformats.org/wordprocessingml/2006/main"

# ---------------------------------------------------------------------
# 1. Insert a new "Meta description" paragraph right after the H1 title
#    paragraph ("Play Clover Lady Free Today! Review & Ratings").
# ---------------------------------------------------------------------
$titlePara = $d.Paragraphs(1)
$titleStart = $titlePara.Range.Start
$titleEnd = $titlePara.Range.End
$titleRange = $d.Range($titleStart, $titleEnd)

$metaXml = @"
<w:p xmlns:w="$w"><w:pPr><w:pStyle w:val="Heading1"/></w:pPr><w:r><w:t>Play Clover Lady Free Today! Review &amp; Ratings</w:t></w:r></w:p><w:p xmlns:w="$w"><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Meta description</w:t></w:r><w:r><w:t>: Explore the enchanting forest with Clover Lady. Read the review, play for free, and discover bonus features, graphics, and design. Compatible on all devices.</w:t></w:r></w:p>
"@
$titleRange.InsertXML($metaXml) | Out-Null

# ---------------------------------------------------------------------
# 2. At the bottom of the document, remove the duplicate title paragraph
#    ("Play Clover Lady Free Today! Review & Ratings") and rewrite the
#    italic paragraph's text to the new feature-image prompt.
# ---------------------------------------------------------------------
$count = $d.Paragraphs.Count
$lastPara = $d.Paragraphs($count)
$prevPara = $d.Paragraphs($count - 1)

$bottomStart = $prevPara.Range.Start
$bottomEnd = $lastPara.Range.End
$bottomRange = $d.Range($bottomStart, $bottomEnd)

$featureText = "For the feature image of Clover Lady, let's have a cartoon-style Maya warrior with glasses. The image should feature the Maya warrior happily playing the game on a mobile device or computer, with the magical forest and mushroom-shaped game grid in the background. The warrior should be holding a clover symbol, with the Metalwolf and girl bonus symbols also visible. The overall style should be colorful and playful, capturing the fairy tale theme of the game."

$featureXml = @"
<w:p xmlns:w="$w"><w:r/><w:r><w:rPr><w:i/></w:rPr><w:t>$featureText</w:t></w:r></w:p>
"@
$bottomRange.InsertXML($featureXml) | Out-Null
